{"js": "// Update the date header and the 3-digit x 1-digit multiplication\n// problems in the practice-sheet table. Each \"old\" value below occurs\n// exactly once in the document, so a body-wide search/replace keyed on\n// the exact old text is unambiguous.\nconst replacements = [\n  [\"2025-04-11 Friday\", \"2025-04-12 Saturday\"],\n  [\"130\\u00d72=\", \"684\\u00d78=\"],\n  [\"978\\u00d76=\", \"149\\u00d72=\"],\n  [\"282\\u00d78=\", \"848\\u00d75=\"],\n  [\"501\\u00d75=\", \"824\\u00d78=\"],\n  [\"340\\u00d73=\", \"217\\u00d77=\"],\n  [\"429\\u00d75=\", \"488\\u00d77=\"],\n  [\"854\\u00d77=\", \"847\\u00d78=\"],\n  [\"988\\u00d72=\", \"901\\u00d77=\"],\n  [\"884\\u00d76=\", \"775\\u00d79=\"],\n  [\"676\\u00d75=\", \"462\\u00d77=\"],\n  [\"952\\u00d77=\", \"588\\u00d73=\"],\n  [\"420\\u00d79=\", \"238\\u00d76=\"],\n  [\"959\\u00d75=\", \"514\\u00d73=\"],\n  [\"764\\u00d77=\", \"745\\u00d75=\"],\n  [\"414\\u00d75=\", \"650\\u00d78=\"],\n  [\"574\\u00d79=\", \"816\\u00d78=\"],\n  [\"527\\u00d72=\", \"812\\u00d72=\"],\n  [\"599\\u00d78=\", \"609\\u00d72=\"],\n  [\"526\\u00d72=\", \"259\\u00d72=\"],\n  [\"852\\u00d77=\", \"746\\u00d72=\"],\n  [\"746\\u00d76=\", \"129\\u00d73=\"],\n  [\"509\\u00d79=\", \"872\\u00d77=\"],\n  [\"705\\u00d79=\", \"954\\u00d74=\"],\n  [\"829\\u00d72=\", \"921\\u00d76=\"],\n  [\"151\\u00d77=\", \"803\\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${oldText}\", found ${found.items.length}`\n    );\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date header and the 3-digit x 1-digit multiplication\n# problems in the practice-sheet table. Each \"old\" value occurs exactly\n# once in the document, so Find/Replace (wdReplaceAll, scoped to the full\n# body) is unambiguous for each pair.\n$d = $word.ActiveDocument\n\nfunction Replace-Unique($oldText, $newText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Find.Execute could not locate expected text: $oldText\"\n    }\n}\n\nReplace-Unique \"2025-04-11 Friday\" \"2025-04-12 Saturday\"\nReplace-Unique \"130\u00d72=\" \"684\u00d78=\"\nReplace-Unique \"978\u00d76=\" \"149\u00d72=\"\nReplace-Unique \"282\u00d78=\" \"848\u00d75=\"\nReplace-Unique \"501\u00d75=\" \"824\u00d78=\"\nReplace-Unique \"340\u00d73=\" \"217\u00d77=\"\nReplace-Unique \"429\u00d75=\" \"488\u00d77=\"\nReplace-Unique \"854\u00d77=\" \"847\u00d78=\"\nReplace-Unique \"988\u00d72=\" \"901\u00d77=\"\nReplace-Unique \"884\u00d76=\" \"775\u00d79=\"\nReplace-Unique \"676\u00d75=\" \"462\u00d77=\"\nReplace-Unique \"952\u00d77=\" \"588\u00d73=\"\nReplace-Unique \"420\u00d79=\" \"238\u00d76=\"\nReplace-Unique \"959\u00d75=\" \"514\u00d73=\"\nReplace-Unique \"764\u00d77=\" \"745\u00d75=\"\nReplace-Unique \"414\u00d75=\" \"650\u00d78=\"\nReplace-Unique \"574\u00d79=\" \"816\u00d78=\"\nReplace-Unique \"527\u00d72=\" \"812\u00d72=\"\nReplace-Unique \"599\u00d78=\" \"609\u00d72=\"\nReplace-Unique \"526\u00d72=\" \"259\u00d72=\"\nReplace-Unique \"852\u00d77=\" \"746\u00d72=\"\nReplace-Unique \"746\u00d76=\" \"129\u00d73=\"\nReplace-Unique \"509\u00d79=\" \"872\u00d77=\"\nReplace-Unique \"705\u00d79=\" \"954\u00d74=\"\nReplace-Unique \"829\u00d72=\" \"921\u00d76=\"\nReplace-Unique \"151\u00d77=\" \"803\u00d78=\"\n"}
